# draft-gandhi-ippm-simple-direct-loss-00.pptx -- "Add files via upload"
#
# Applies the four content edits described by the commit's diff:
#   1. Handout-master footer date field text: 2/24/21 -> 2/28/21
#   2. Slide 3 ("Requirements and Scope") bullet text tweak
#   3. Slide 7 table header cell: "...TLV?" -> "...TLV2?"
#   4. Slide 9 ("Next Steps") new italic line "Define New STAMP Direct
#      Measurement TLV2?" added as a new run in the trailing empty paragraph
#
# (The p14:modId GUID-ish counter on slide 7's table, and the namespace
# attribute re-ordering on the Mac placeholderFlag extension in
# slideLayout12.xml, are opaque/internal bookkeeping values with no
# corresponding PowerPoint object model surface -- they are not
# user content and are left to the host application to manage.)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Handout Master date placeholder: "2/24/21" -> "2/28/21"
# ---------------------------------------------------------------------
$hm = $p.HandoutMaster
$dtf = $hm.HeadersFooters.DateAndTime
$dtf.UseFormat = $false
$dtf.Value = "2/28/21"

# ---------------------------------------------------------------------
# 2) Slide 3 - "Requirements and Scope": update the "Goals" bullet
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$content3 = $s3.Shapes.Item(2)
$tr3 = $content3.TextFrame.TextRange
$goalPara = $tr3.Paragraphs(7)
# Replace only the exact run span (Start/Length) so the untouched text on
# either side is not re-diffed into extra runs.
$goalRange = $tr3.Characters($goalPara.Start, $goalPara.Length)
$goalRange.Text = "Avoid provisioning and maintaining each test session on Session-Reflector"

# ---------------------------------------------------------------------
# 3) Slide 7 - table cell "Case 1a. ... TLV?" -> "...TLV2?"
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$tableShape = $s7.Shapes.Item(2)
$tbl = $tableShape.Table
$cell = $tbl.Cell(1, 3)
$cell.Shape.TextFrame.TextRange.Text = "Case 1a. Define New STAMP Direct Measurement TLV2?"

# ---------------------------------------------------------------------
# 4) Slide 9 - "Next Steps": add new italic line in the trailing empty
#    paragraph of the content placeholder.
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$content9 = $s9.Shapes.Item(2)
$tr9 = $content9.TextFrame.TextRange
$lastPara = $tr9.Paragraphs(3)
$lastPara.Text = "Define New STAMP Direct Measurement TLV2?"

# Re-fetch and apply the run formatting used elsewhere on this slide
# (24pt, italic) to match the new line's intended look.
$newPara = $s9.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3)
$newPara.Font.Size = 24
$newPara.Font.Italic = $true
